$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "29.883.41"
$ws.Range("E2").Value = "  +0.10%  "
$ws.Range("D3").Value = "1.893.58"
$ws.Range("E3").Value = "  -0.05%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "0.7720"
$ws.Range("E5").Value = "  -2.30%  "
$ws.Range("D6").Value = "244.39"
$ws.Range("E6").Value = "  +0.26%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "0.3124"
$ws.Range("E8").Value = "  -0.80%  "
$ws.Range("D9").Value = "25.71"
$ws.Range("E9").Value = "  +1.56%  "
$ws.Range("D10").Value = "0.07233"
$ws.Range("E10").Value = "  -0.18%  "
$ws.Range("D11").Value = "0.08722"
$ws.Range("E11").Value = "  +7.81%  "
$ws.Range("D12").Value = "2.080.70"
$ws.Range("E12").Value = "  +8.24%  "
$ws.Range("D13").Value = "0.7717"
$ws.Range("E13").Value = "  +0.71%  "
$ws.Range("D14").Value = "5.426"
$ws.Range("E14").Value = "  -2.11%  "
$ws.Range("D15").Value = "94.31"
$ws.Range("E15").Value = "  +2.00%  "
$ws.Range("D16").Value = "6.213"
$ws.Range("E16").Value = "  +0.79%  "
$ws.Range("D17").Value = "30.019.78"
$ws.Range("E17").Value = "  +0.48%  "
$ws.Range("D18").Value = "13.93"
$ws.Range("E18").Value = "  +0.11%  "
$ws.Range("D19").Value = "245.11"
$ws.Range("E19").Value = "  +0.49%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "0.000007878"
$ws.Range("E20").Value = "  +1.39%  "
$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").Value = "2.263.32"
$ws.Range("E21").Value = "  +4.57%  "
$ws.Range("D22").Value = "8.185"
$ws.Range("E22").Value = "  +0.47%  "
$ws.Range("D23").Value = "1.002"
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("E25").Value = "  -2.74%  "
$ws.Range("D26").Value = "9.522"
$ws.Range("E26").Value = "  +1.21%  "
$ws.Range("D27").Value = "162.32"
$ws.Range("E27").Value = "  -0.61%  "
$ws.Range("D28").Value = "18.80"
$ws.Range("E28").Value = "  +0.37%  "
$ws.Range("D29").Value = "2.039"
$ws.Range("E29").Value = "  -0.84%  "
$ws.Range("D30").Value = "1.429"
$ws.Range("E30").Value = "  +1.97%  "
$ws.Range("D31").Value = "1.543"
$ws.Range("E31").Value = "  -0.10%  "
$ws.Range("D32").Value = "4.534"
$ws.Range("E32").Value = "  +0.78%  "
$ws.Range("D33").Value = "4.121"
$ws.Range("E33").Value = "  +0.32%  "
$ws.Range("D34").Value = "0.05484"
$ws.Range("E34").Value = "  -1.57%  "
$ws.Range("E35").Value = "  -1.68%  "
$ws.Range("D36").Value = "0.7543"
$ws.Range("E36").Value = "  +1.71%  "
$ws.Range("D37").Value = "1.003"
$ws.Range("E37").Value = "  +0.03%  "
$ws.Range("D38").Value = "2.720"
$ws.Range("E38").Value = "  +3.84%  "
$ws.Range("D39").Value = "0.01963"
$ws.Range("E39").Value = "  +1.90%  "
$ws.Range("E40").Value = "  +0.42%  "
$ws.Range("D41").Value = "0.4505"
$ws.Range("E41").Value = "  +1.71%  "
$ws.Range("D42").Value = "73.62"
$ws.Range("E42").Value = "  -0.57%  "
$ws.Range("D43").Value = "1.094.46"
$ws.Range("E43").Value = "  -4.06%  "
$ws.Range("D44").Value = "6.040"
$ws.Range("E44").Value = "  +2.71%  "
$ws.Range("D45").Value = "0.8562"
$ws.Range("E45").Value = "  +0.67%  "
$ws.Range("D46").Value = "1.001"
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").Value = "103.03"
$ws.Range("E47").Value = "  -1.40%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "1.885"
$ws.Range("E48").Value = "  +0.43%  "
$ws.Range("B49").Value = "Aptos"
$ws.Range("C49").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D49").Value = "7.625"
$ws.Range("E49").Value = "  +2.28%  "
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "2.133.88"
$ws.Range("E50").Value = "  +2.67%  "
$ws.Range("D51").Value = "9.841"
$ws.Range("E51").Value = "  -1.55%  "
